$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 389, shifting existing rows 389..426 down to 390..427.
$ws.Rows.Item(389).Insert()

# Populate the newly inserted row 389 with the new record.
$ws.Cells.Item(389, 1).Value = 5
$ws.Cells.Item(389, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(389, 3).Value = "Maule"
$ws.Cells.Item(389, 4).Value = 45132
$ws.Cells.Item(389, 5).Value = 7
$ws.Cells.Item(389, 6).Value = "Fruta"
$ws.Cells.Item(389, 7).Value = 100108
$ws.Cells.Item(389, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(389, 9).Value = 100108005
$ws.Cells.Item(389, 10).Value = "Piña"
$ws.Cells.Item(389, 11).Value = "Caramelo"
$ws.Cells.Item(389, 12).Value = "Segunda"
$ws.Cells.Item(389, 13).Value = 120
$ws.Cells.Item(389, 14).Value = 19000
$ws.Cells.Item(389, 15).Value = 19000
$ws.Cells.Item(389, 16).Value = 19000
$ws.Cells.Item(389, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(389, 18).Value = "Ecuador"
$ws.Cells.Item(389, 19).Value = 1357
$ws.Cells.Item(389, 20).Value = 14

Write-Output "done"
